# Fix arrondie <etat virement>
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (A:M) for rows 2-9 (detail lines) and row 10 (totals),
# replacing the previous 5 detail rows + 1 total row (old rows 2-7).
$data = @(
    @("910/TANGER /AV1",     "Direction régionale", "D235689",  "KAMILIA LALA", "non", "mensuelle", 15, 24000,    24000,    3600,     3600,     "--", 20400),
    @("910/TANGER /AV1",     "Direction régionale", "K3544354", "ABDOU FAFA",   "non", "mensuelle", 15, 13354.4,  13354.4,  2003.16,  2003.16,  "--", 11351.24),
    @("910/TANGER /AV1",     "Direction régionale", "L254654",  "SAMIR DADA",   "non", "mensuelle", 10, 2645.6,   2645.6,   264.56,   264.56,   "--", 2381.04),
    @("115/TANGER MED/AV1",  "Point de vente",      "L5245475", "MORAD JOJO",   "non", "mensuelle", 10, 9000,     0,        900,      0,        "--", 8100),
    @("115/TANGER MED/AV1",  "Point de vente",      "Z213568",  "NABIL MOMO",   "non", "mensuelle", 10, 6000,     0,        600,      0,        "--", 5400),
    @("115/TANGER MED/AV1",  "Point de vente",      "L525655",  "KHALID RARA",  "non", "mensuelle", 10, 6000,     0,        600,      0,        "--", 5400),
    @("844/T-SUD",           "Point de vente",      "K324554",  "KARIMA SASA",  "non", "mensuelle", 15, 24000,    0,        3600,     0,        "--", 20400),
    @("844/T-SUD",           "Point de vente",      "IL12254",  "FARIDA VAVA",  "non", "mensuelle", 10, 6000,     0,        600,      0,        "--", 5400)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value  = $row[4]
    $ws.Cells.Item($r, 6).Value  = $row[5]
    $ws.Cells.Item($r, 7).Value  = $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
}

# Row 10: totals row, shifted down from the old row 7 with updated values.
$ws.Cells.Item(10, 1).Value  = " "
$ws.Cells.Item(10, 2).Value  = " "
$ws.Cells.Item(10, 3).Value  = " "
$ws.Cells.Item(10, 4).Value  = " "
$ws.Cells.Item(10, 5).Value  = " "
$ws.Cells.Item(10, 6).Value  = " "
$ws.Cells.Item(10, 7).Value  = " "
$ws.Cells.Item(10, 8).Value  = 91000
$ws.Cells.Item(10, 9).Value  = 0
$ws.Cells.Item(10, 10).Value = 12167.72
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 78832.28
